$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 3088.25
$ws.Range("J17").Value = 3088.25
$ws.Range("L17").Value = 9264.75
$ws.Range("N17").Value = -9600.75
# Row 34
$ws.Range("H34").Value = 3874.75
$ws.Range("I34").Value = 3499.6667
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 3499.6667
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -3296.6667
$ws.Range("N34").Value = -5406
# Row 36
$ws.Range("H36").Value = 3874.75
$ws.Range("I36").Value = 3499.6667
$ws.Range("J36").Value = 5000
$ws.Range("K36").Value = 3499.6667
$ws.Range("L36").Value = 5000
$ws.Range("M36").Value = -2784.6667
$ws.Range("N36").Value = -6430
# Row 43
$ws.Range("H43").Value = 999.6667
$ws.Range("J43").Value = 999
$ws.Range("L43").Value = 999
$ws.Range("N43").Value = -1137
# Row 98
$ws.Range("H98").Value = 1464.6923
$ws.Range("I98").Value = 1420.9166
$ws.Range("K98").Value = 1420.9166
$ws.Range("M98").Value = 77.08339999999998
# Row 103
$ws.Range("H103").Value = 1782.7142
$ws.Range("I103").Value = 10004
$ws.Range("J103").Value = 412.5
$ws.Range("K103").Value = 30012
$ws.Range("L103").Value = 1237.5
$ws.Range("M103").Value = -29426
$ws.Range("N103").Value = -2409.5
# Row 107
$ws.Range("H107").Value = 477.53845
$ws.Range("I107").Value = 382.72726
$ws.Range("K107").Value = 382.72726
$ws.Range("M107").Value = 1537.27274
# Row 122
$ws.Range("H122").Value = 1464.6923
$ws.Range("I122").Value = 1420.9166
$ws.Range("K122").Value = 4262.7498
$ws.Range("M122").Value = -1812.7498

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12089.782
$ws.Range("I32").Value = 9836.223
$ws.Range("K32").Value = 9836.223
$ws.Range("M32").Value = -9549.223
# Row 61
$ws.Range("H61").Value = 6249
$ws.Range("I61").Value = 6283.75
$ws.Range("K61").Value = 6283.75
$ws.Range("M61").Value = -6071.75
# Row 122
$ws.Range("H122").Value = 2999.5
$ws.Range("I122").Value = 2999
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 8997
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -6547
$ws.Range("N122").Value = -13900
# Row 132
$ws.Range("H132").Value = 1476.1698
$ws.Range("I132").Value = 1489.1923
$ws.Range("K132").Value = 4467.5769
$ws.Range("M132").Value = -1937.5769
# Row 136
$ws.Range("H136").Value = 6249
$ws.Range("I136").Value = 6283.75
$ws.Range("K136").Value = 18851.25
$ws.Range("M136").Value = -16301.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 31333.334
$ws.Range("I82").Value = 6000
$ws.Range("K82").Value = 6000
$ws.Range("M82").Value = -5617
# Row 85
$ws.Range("H85").Value = 31333.334
$ws.Range("I85").Value = 6000
$ws.Range("K85").Value = 6000
$ws.Range("M85").Value = -4674

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5954.0435
$ws.Range("I31").Value = 5535.3335
$ws.Range("K31").Value = 5535.3335
$ws.Range("M31").Value = -5240.3335
# Row 34
$ws.Range("H34").Value = 5954.0435
$ws.Range("I34").Value = 5535.3335
$ws.Range("K34").Value = 5535.3335
$ws.Range("M34").Value = -5333.3335
# Row 60
$ws.Range("H60").Value = 19624.777
$ws.Range("J60").Value = 19874.715
$ws.Range("L60").Value = 19874.715
$ws.Range("N60").Value = -20896.715
# Row 62
$ws.Range("H62").Value = 4998.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4998.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4998.5
$ws.Range("N62").Value = -6246.5
$ws.Range("M62").ClearContents()
# Row 65
$ws.Range("H65").Value = 4998.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4998.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 24992.5
$ws.Range("N65").Value = -31232.5
$ws.Range("M65").ClearContents()
# Row 122
$ws.Range("H122").Value = 3024.52
$ws.Range("I122").Value = 2860.5
$ws.Range("J122").Value = 3446.2856
$ws.Range("K122").Value = 8581.5
$ws.Range("L122").Value = 10338.8568
$ws.Range("M122").Value = -6131.5
$ws.Range("N122").Value = -15238.8568
# Row 132
$ws.Range("H132").Value = 4042.75
$ws.Range("I132").Value = 2425.4285
$ws.Range("K132").Value = 7276.2855
$ws.Range("M132").Value = -4746.2855
# Row 141
$ws.Range("H141").Value = 207245.1
$ws.Range("J141").Value = 207245.1
$ws.Range("L141").Value = 207245.1
$ws.Range("N141").Value = -217605.1

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
# Row 87
$ws.Range("H87").Value = 8000.091
$ws.Range("I87").Value = 6300.1
$ws.Range("K87").Value = 18900.3
$ws.Range("M87").Value = -17652.3
# Row 90
$ws.Range("H90").Value = 8000.091
$ws.Range("I90").Value = 6300.1
$ws.Range("K90").Value = 56700.9
$ws.Range("M90").Value = -50460.9
# Row 134
$ws.Range("H134").Value = 993.619
$ws.Range("I134").Value = 943.3
$ws.Range("K134").Value = 2829.9
$ws.Range("M134").Value = 2240.1
# Row 139
$ws.Range("H139").Value = 8570.632
$ws.Range("I139").Value = 3788.1667
$ws.Range("K139").Value = 11364.5001
$ws.Range("M139").Value = -6224.500100000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 14
$ws.Range("H14").Value = 4286.75
$ws.Range("I14").Value = 2349.5
$ws.Range("J14").Value = 4932.5
$ws.Range("K14").Value = 2349.5
$ws.Range("L14").Value = 4932.5
$ws.Range("M14").Value = -2181.5
$ws.Range("N14").Value = -5268.5
# Row 80
$ws.Range("H80").Value = 2974.3333
$ws.Range("I80").Value = 2056.8
$ws.Range("J80").Value = 3433.1
$ws.Range("K80").Value = 2056.8
$ws.Range("L80").Value = 3433.1
$ws.Range("M80").Value = -1058.8
$ws.Range("N80").Value = -5429.1
# Row 83
$ws.Range("H83").Value = 2974.3333
$ws.Range("I83").Value = 2056.8
$ws.Range("J83").Value = 3433.1
$ws.Range("K83").Value = 10284
$ws.Range("L83").Value = 17165.5
$ws.Range("M83").Value = -5292
$ws.Range("N83").Value = -27149.5
# Row 102
$ws.Range("H102").Value = 2957.9333
$ws.Range("I102").Value = 1549.04
$ws.Range("K102").Value = 1549.04
$ws.Range("M102").Value = 72.96000000000004

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 38
$ws.Range("H38").Value = 350000
$ws.Range("I38").Value = 510000
$ws.Range("J38").Value = 30000
$ws.Range("K38").Value = 510000
$ws.Range("L38").Value = 30000
$ws.Range("M38").Value = -509590
$ws.Range("N38").Value = -30820
# Row 122
$ws.Range("H122").Value = 2621.5557
$ws.Range("I122").Value = 2599.25
$ws.Range("K122").Value = 7797.75
$ws.Range("M122").Value = -5347.75
# Row 132
$ws.Range("H132").Value = 8986.218999999999
$ws.Range("I132").Value = 9115.966
$ws.Range("K132").Value = 27347.898
$ws.Range("M132").Value = -24817.898

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
# Row 54
$ws.Range("H54").Value = 47904.89
$ws.Range("J54").Value = 99015.336
$ws.Range("L54").Value = 99015.336
$ws.Range("N54").Value = -100055.336
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("M95").ClearContents()
$ws.Range("N95").ClearContents()
